$d = $word.ActiveDocument
$xml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>הפתרון האופטימלי הוא</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> ליצור כמה </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>מיקרוסרוויסים</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> בשביל שרתי ה</w:t>
      </w:r>
      <w:r>
        <w:t>backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> ולהשתמש ב</w:t>
      </w:r>
      <w:r>
        <w:t>Load Balancer</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> (</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>LB</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">) על מנת ליצור שוויון בחלוקת העומס. יצירת כמה שרתי </w:t>
      </w:r>
      <w:r>
        <w:t>backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> תעזור בהפחתת עומס המשתמשים על שרת אחד, ובכך תעזור בהימנעות מקריסת האתר עקב עומס.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> כאשר יש כמה שרתי </w:t>
      </w:r>
      <w:r>
        <w:t>backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">, גם אם אחד מהשרתים קורס, הוא לא </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>יקריס</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> את כל האתר, מכיוון שה-</w:t>
      </w:r>
      <w:r>
        <w:t>LB</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> יחוון את המשתמשים לשרת </w:t>
      </w:r>
      <w:r>
        <w:t>backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> אחר פעיל.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:color w:val="FF0000"/>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:color w:val="FF0000"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">מכיוון שאין הרבה זמן, וגם העומס על האתר שאבנה לא אמור להיות כבד, לא אעשה זאת באתר, ויהיה רק שרת </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:color w:val="FF0000"/>
        </w:rPr>
        <w:t>backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:color w:val="FF0000"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> אחד.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>H</w:t>
      </w:r>
      <w:r>
        <w:t>orizontal Scaling</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> יצירת עוד </w:t>
      </w:r>
      <w:r>
        <w:t>nodes</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> כדי לחלק את העומס (למשל, יצירת עוד שרתי </w:t>
      </w:r>
      <w:r>
        <w:t>backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>)</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:t>Vertical Scaling</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> הוספת משאבים לשרת אחד. למשל, להגדיל את כוח ההתמודדות של השרת עם עומס.</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> יותר קל לעשות אבל פחות חסכוני. בנוסף יש סיכון של </w:t>
      </w:r>
      <w:r>
        <w:t>Single Point Of Failure</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> כי יש רק שרת אחד עם כוח גדול, ואם הוא נופל כל האתר קורס.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>אני אבחר להשתמש ב</w:t>
      </w:r>
      <w:r>
        <w:t>node.js</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> על מנת לבנות את שרת ה</w:t>
      </w:r>
      <w:r>
        <w:t>backend</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> שלי כי יש לו יתרונות שיעזרו לי לעשות את האתר יותר </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>סקלבילי</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> כמו</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>: התמודדות א-סינכרונית עם בקשות כך שניתן לקבל הרבה בקשות בו-זמנית ללא בעיה.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>כדי להימנע מ</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>SPOF</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> גם בשכבת ה</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>DB</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>, הפתרון האידיאלי יהיה להשתמש ב-</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>D</w:t>
      </w:r>
      <w:r>
        <w:t>ata Replication</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">. יהיה </w:t>
      </w:r>
      <w:r>
        <w:t>master DB</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> ועותקים של המאסטר. במאסטר יהיה ניתן רק לכתוב, ואת העותקים שלו יהיה ניתן רק לקרוא. יתרונות השיטה: המודל מאפשר ליותר שאילתות להתבצע בו-זמנית ואם ה</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>DB</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> נהרס פיזית, המידע נשמר, מכיוון שיש לו עותקים במקומות אחרים.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>על מנת לשפר את זמן התגובה והטעינה:</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">האידיאל </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> להשתמש ב</w:t>
      </w:r>
      <w:r>
        <w:t>caching</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> ו</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>CDN</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>C</w:t>
      </w:r>
      <w:r>
        <w:t>aching</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> מקום אחסון זמני ל</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">תוצאות של תגובות יקרות או של מידע שניגשים אליו בתדירות גבוהה, כדי שלבקשות בשבילן ייקח הרבה פחות זמן להחזיר תשובה ולהגיב. </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>כאשר אנו מקבלים בקשה, שרת האינטרנט בודק אם המידע נמצא ב</w:t>
      </w:r>
      <w:r>
        <w:t>cache</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>, במידה וכן, הוא מחזיר את המידע משם, ובמידה ולא, הוא שולח בקשה ל</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>DB</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>, שומר את התשובה ב</w:t>
      </w:r>
      <w:r>
        <w:t>cache</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>, ומחזיר את המידע ללקוח.</w:t>
      </w:r>
    </w:p>
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>CDN</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rtl/>
        </w:rPr>
        <w:t>–</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> רשת שרתים מרוחקים שמשומשים על מנת לאחסן מידע סטטי, כמו תמונות, סרטונים, קבצי </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>CSS</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve">, קבצי </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>JS</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t xml:space="preserve"> ועוד. כאשר הלקוח מבקש מידע סטטי, השרת האינטרנטי בודק אם הוא קיים ב</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>CDN</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>, במידה וכן, שולח ללקוח, ובמידה ולא, מבקש מהשרת, מאחסן את התשובה ב</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
        </w:rPr>
        <w:t>CDN</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="cs"/>
          <w:rtl/>
        </w:rPr>
        <w:t>, ומחזיר ללקוח.</w:t>
      </w:r>
    </w:p>
'@
$target = $d.Paragraphs(2).Range
$target.InsertXML($xml)
Write-Output ("ParaCount=" + $d.Paragraphs.Count)
